# Venue slide: bump the meeting month from May to June.
#
# Slide 4 ("May meetup venue") has:
#   - a Title shape whose whole text is "May meetup venue"
#   - a body text box whose first paragraph ends in
#     "...kindly hosting us for our May meetup" (followed by two more
#     paragraphs - Location / Website - that must stay untouched)

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(4)

# Title 1: "May meetup venue" -> "June meetup venue"
$title = $s.Shapes.Item(1)
$title.TextFrame.TextRange.Text = "June meetup venue"

# Content Placeholder 2 (the shape holding the venue blurb): update just
# the first paragraph/run so the Location/Website paragraphs (with their
# superscript "th" run) are left completely untouched.
$body = $s.Shapes.Item(7)
$firstPara = $body.TextFrame.TextRange.Paragraphs(1)
$firstPara.Runs(1).Text = "Incremental Group are kindly hosting us for our June meetup"
